# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh case numbers for several countries (unchanged A/D/F columns where applicable) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1004645
$ws.Range("C4").Value = 17485
$ws.Range("E4").Value = 810685
$ws.Range("G4").Value = 1082
$ws.Range("H4").Value = 56495

# Row 8: Alemania
$ws.Range("B8").Value = 158389
$ws.Range("C8").Value = 619
$ws.Range("E8").Value = 37839
$ws.Range("G8").Value = 74
$ws.Range("H8").Value = 6050

# Row 14: Brasil
$ws.Range("B14").Value = 66501
$ws.Range("C14").Value = 3642
$ws.Range("E14").Value = 31806
$ws.Range("G14").Value = 272
$ws.Range("H14").Value = 4543

# Row 22: Ecuador
$ws.Range("B22").Value = 23240
$ws.Range("C22").Value = 521
$ws.Range("E22").Value = 21211
$ws.Range("G22").Value = 87
$ws.Range("H22").Value = 663

# Row 26: Israel
$ws.Range("B26").Value = 15555
$ws.Range("C26").Value = 112
$ws.Range("D26").Value = 7200
$ws.Range("E26").Value = 8151
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 204

# --- Swap Costa Rica / Niger rows + update Niger's figures ---
# Before: row 99 = Costa Rica (697,2,287,404,8,0,6); row 100 = Niger (696,0,350,317,0,0,29)
# After:  row 99 = Niger  (701,5,385,287,0,0,29);   row 100 = Costa Rica (697,2,287,404,8,0,6)
$ws.Range("A99").Value = "Niger"
$ws.Range("B99").Value = 701
$ws.Range("C99").Value = 5
$ws.Range("D99").Value = 385
$ws.Range("E99").Value = 287
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 29

$ws.Range("A100").Value = "Costa Rica"
$ws.Range("B100").Value = 697
$ws.Range("C100").Value = 2
$ws.Range("D100").Value = 287
$ws.Range("E100").Value = 404
$ws.Range("F100").Value = 8
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 6

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Abril de 2020 a las 22:22"
